$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D values are plain numeric-looking strings that must remain
# TEXT (as in the source file, which stores them as inline strings).
# Excel's COM layer auto-converts numeric-looking input to real numbers,
# so we briefly mark each cell as Text before assigning the value, then
# clear the formatting again so the cell keeps the workbook's original
# (default) style -- only its text content changes.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.834.07'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.624.47'
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.08'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.17'
$ws.Range("D6").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.150'
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.66'
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.097.07'
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.689.64'
$ws.Range("D15").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.603.08'
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.31'
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.63'
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.94'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.88'
$ws.Range("D21").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.32'
$ws.Range("D24").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.19'
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '562.07'
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.20'
$ws.Range("D29").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0841'
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.73'
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.21'
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '169.33'
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.408'
$ws.Range("D37").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.94'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.36'
$ws.Range("D40").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '170.10'
$ws.Range("D42").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.92'
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.37'
$ws.Range("D46").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0247'
$ws.Range("D48").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0967'
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.23'
$ws.Range("D51").ClearFormats()

# Column E values (percentage change) already contain non-numeric
# characters (%, padding spaces) so they remain text automatically.
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("E3").Value = '  -1.45%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("E11").Value = '  +3.09%  '
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("E14").Value = '  -1.24%  '
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("E18").Value = '  +6.84%  '
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("E25").Value = '  +12.03%  '
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("E29").Value = '  +4.38%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("E44").Value = '  +3.12%  '
$ws.Range("E45").Value = '  +3.43%  '
$ws.Range("E46").Value = '  -5.32%  '
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("E49").Value = '  +4.68%  '
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("E51").Value = '  +1.72%  '
